$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.099.85"
$ws.Range("E2").Value = "  -4.31%  "
$ws.Range("D3").Value = "2.903.59"
$ws.Range("E3").Value = "  -3.55%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.64"
$ws.Range("E5").Value = "  -5.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.99"
$ws.Range("E6").Value = "  -7.25%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -2.23%  "
$ws.Range("D9").Value = "2.904.60"
$ws.Range("E9").Value = "  -3.79%  "
$ws.Range("E10").Value = "  -5.40%  "
$ws.Range("E11").Value = "  -8.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.354"
$ws.Range("E12").Value = "  -2.93%  "
$ws.Range("D13").Value = "3.410.34"
$ws.Range("E13").Value = "  -3.42%  "
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "60.266.16"
$ws.Range("E15").Value = "  -4.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.62"
$ws.Range("E16").Value = "  -5.47%  "
$ws.Range("D17").Value = "2.892.55"
$ws.Range("E17").Value = "  -4.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000139"
$ws.Range("E18").Value = "  -6.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.92"
$ws.Range("E19").Value = "  -3.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.48"
$ws.Range("E20").Value = "  -3.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "359.82"
$ws.Range("E21").Value = "  -9.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.55"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.21"
$ws.Range("E25").Value = "  -3.02%  "
$ws.Range("D26").Value = "3.020.70"
$ws.Range("E26").Value = "  -3.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.448"
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.172"
$ws.Range("E28").Value = "  -8.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("E30").Value = "  -10.75%  "
$ws.Range("E31").Value = "  -12.52%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  -5.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.40"
$ws.Range("E34").Value = "  -5.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "151.04"
$ws.Range("E35").Value = "  -6.01%  "
$ws.Range("E36").Value = "  -8.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.52"
$ws.Range("E37").Value = "  -8.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.985"
$ws.Range("E38").Value = "  -9.99%  "
$ws.Range("E39").Value = "  -8.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.67"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "2.337.31"
$ws.Range("E41").Value = "  -6.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.45"
$ws.Range("E42").Value = "  -8.09%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.65"
$ws.Range("E43").Value = "  -6.64%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.643"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.66"
$ws.Range("E45").Value = "  -8.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0567"
$ws.Range("E46").Value = "  -4.68%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  -4.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.36"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("E50").Value = "  -6.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0926"
$ws.Range("E51").Value = "  -2.08%  "
